# Natmi following Dr Hou advice
# Update Ligand-expressing cells (E) and Receptor-expressing cells (K) counts
# from 1 to 3, and recompute the dependent expression/specificity columns
# (G, H, M, N, O, P, Q, R, S, T) for data rows 2-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ E = 3; G = 2.822099333333334; H = 8.466298; K = 3; M = 9.403031; N = 28.209093; O = 0.1079186507462474; P = 0.1079186507462474; Q = 26.53628751641267; R = 238.826587647714; S = 0.1079186507462474; T = 0.1079186507462474 }
    3  = @{ E = 3; G = 2.822099333333334; H = 8.466298; K = 3; M = 21.85073566666667; N = 65.552207; O = 0.2507810418746436; P = 0.2507810418746436; Q = 61.66494655774289; R = 554.984519019686; S = 0.2507810418746436; T = 0.2507810418746436 }
    4  = @{ E = 3; G = 2.822099333333334; H = 8.466298; K = 3; M = 14.75061866666667; N = 44.251856; O = 0.1692929507707757; P = 0.1692929507707757; Q = 41.62771110545423; R = 374.649399949088; S = 0.1692929507707757; T = 0.1692929507707757 }
    5  = @{ E = 3; G = 2.822099333333334; H = 8.466298; K = 3; M = 19.24455733333333; N = 57.733672; O = 0.2208699154157989; P = 0.2208699154157988; Q = 54.31005242069511; R = 488.790471786256; S = 0.2208699154157989; T = 0.2208699154157988 }
    6  = @{ E = 3; G = 2.822099333333334; H = 8.466298; K = 3; M = 6.960233; N = 20.880699; O = 0.07988264148437944; P = 0.07988264148437943; Q = 19.64246890914467; R = 176.782220182302; S = 0.07988264148437944; T = 0.07988264148437943 }
    7  = @{ E = 3; G = 2.822099333333334; H = 8.466298; K = 3; M = 14.921556; N = 44.764668; O = 0.171254799708155; P = 0.171254799708155; Q = 42.110113239896; R = 378.991019159064; S = 0.171254799708155; T = 0.171254799708155 }
}

foreach ($row in $values.Keys) {
    $rowData = $values[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
